$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values could otherwise be
# auto-converted to numbers (losing the literal string representation),
# matching the original inlineStr text cells.
$textCells = @("D5","D6","D9","D10","D12","D15","D19","D21","D22","D23","D24","D25","D27","D28","D30","D31","D32","D33","D35","D37","D38","D39","D40","D41","D42","D44","D45","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "45.876.23"
$ws.Range("E2").Value = "  +7.96%  "
$ws.Range("D3").Value = "2.419.61"
$ws.Range("E3").Value = "  +5.91%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "115.04"
$ws.Range("E5").Value = "  +11.66%  "
$ws.Range("D6").Value = "319.63"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "0.632"
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("D10").Value = "43.19"
$ws.Range("E10").Value = "  +11.72%  "
$ws.Range("E11").Value = "  +4.70%  "
$ws.Range("D12").Value = "8.76"
$ws.Range("E12").Value = "  +6.97%  "
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("E14").Value = "  +4.91%  "
$ws.Range("D15").Value = "16.02"
$ws.Range("E15").Value = "  +5.47%  "
$ws.Range("D16").Value = "2.765.75"
$ws.Range("E16").Value = "  +5.04%  "
$ws.Range("D17").Value = "2.414.87"
$ws.Range("E17").Value = "  +6.08%  "
$ws.Range("D18").Value = "45.835.68"
$ws.Range("E18").Value = "  +7.42%  "
$ws.Range("D19").Value = "7.65"
$ws.Range("E19").Value = "  +5.31%  "
$ws.Range("E20").Value = "  +4.42%  "
$ws.Range("D21").Value = "13.55"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "75.07"
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("D23").Value = "3.54"
$ws.Range("E23").Value = "  +4.79%  "
$ws.Range("D24").Value = "270.27"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  +7.99%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").Value = "7.71"
$ws.Range("E27").Value = "  +8.48%  "
$ws.Range("D28").Value = "11.37"
$ws.Range("E28").Value = "  +5.84%  "
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "39.27"
$ws.Range("E30").Value = "  +10.55%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "23.11"
$ws.Range("E31").Value = "  +3.76%  "
$ws.Range("D32").Value = "0.0973"
$ws.Range("E32").Value = "  +14.64%  "
$ws.Range("D33").Value = "172.68"
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("E34").Value = "  +17.16%  "
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").Value = "  +8.89%  "
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").Value = "4.99"
$ws.Range("E37").Value = "  +10.53%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  +13.76%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "4.14"
$ws.Range("E39").Value = "  +15.02%  "
$ws.Range("D40").Value = "0.0368"
$ws.Range("E40").Value = "  +6.52%  "
$ws.Range("D41").Value = "1.79"
$ws.Range("E41").Value = "  +15.64%  "
$ws.Range("D42").Value = "103.51"
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("E43").Value = "  +6.96%  "
$ws.Range("D44").Value = "72.06"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "13.43"
$ws.Range("E45").Value = "  +11.83%  "
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "5.88"
$ws.Range("E47").Value = "  +14.61%  "
$ws.Range("D48").Value = "118.57"
$ws.Range("E48").Value = "  +7.31%  "
$ws.Range("E49").Value = "  +18.00%  "
$ws.Range("D50").Value = "9.40"
$ws.Range("E50").Value = "  +9.15%  "
$ws.Range("D51").Value = "79.68"
$ws.Range("E51").Value = "  +3.64%  "
